$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the title from A1 to D1, keeping its original formatting
$ws.Range("A1").Copy($ws.Range("D1"))
$ws.Range("A1").Clear()

# Row 4: update points value for "The quality of your code"
$ws.Range("B4").Value = 20

# Row 9: restyle "Completeness of your unit tests" label (bold -> not bold)
$ws.Range("D9").Font.Bold = $false

# Fill in new row 10 for "Completeness of functional testing" (the row
# was previously blank, so no shifting of existing rows is needed).
# Re-use D9's (now un-bolded) formatting rather than setting font
# properties from scratch.
$ws.Range("B10").Value = 10
$ws.Range("D10").Value = "Completeness of functional testing"
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)

# Row 11: update points value and restyle label (bold -> not bold)
$ws.Range("B11").Value = 35
$ws.Range("D11").Font.Bold = $false

# Add a new row 20 with the grand total
$ws.Range("B20").Formula = "=SUM(B4:B18)"
$ws.Range("D20").Value = "Total"
$ws.Range("D20").Font.Bold = $true

# Column width tweaks (values chosen so the saved OOXML width, which is
# quantized to 1/6 character-width steps by this runtime, lands as close
# as possible to the target widths from the diff)
$ws.Columns("A").ColumnWidth = 5.3
$ws.Columns("B").ColumnWidth = 7
$ws.Columns("C").ColumnWidth = 1.6
$ws.Columns("G").ColumnWidth = 9.6

# Selection moves to A3
$ws.Range("A3").Select() | Out-Null
